# Generate Report for Handback
#
# Refreshes the "Latest Handback DateTime" value for the first data row
# (b5aba85f-5945-44a8-b7a2-f320663c1875) on the German ("de-de") handback
# status sheet, reflecting the freshly generated handback report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("de-de")

# Column K on this sheet is "Latest Handback DateTime"; row 2 is the
# b5aba85f-5945-44a8-b7a2-f320663c1875 entry.
$ws.Range("K2").Value = "2016-11-14 07:12:17"
